# Refresh the GSC (Google Search Console) "Video Indexing" export.
#
# The export window rolled forward by one day: the oldest day in the
# "Chart" sheet (2025-11-18, which had no data yet - blank indexed
# counts and 0 impressions) has aged out of the report, so its row is
# removed and every following day's row shifts up to take its place.
# The "Table" and "Metadata" sheets are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest date (2025-11-18). Delete it; Excel shifts
# rows 3:88 up to rows 2:87, so the sheet ends up with one fewer row
# (A1:D87) and every date advances one row closer to the header.
$ws.Range("A2").EntireRow.Delete()
